$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.113.64"
$ws.Range("E2").Value = "  +0.59%  "

$ws.Range("D3").Value = "2.598.43"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'574.91"
$ws.Range("E5").Value = "  +3.21%  "

$ws.Range("D6").Value = "'142.54"
$ws.Range("E6").Value = "  +0.68%  "

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("E8").Value = "  -0.27%  "

$ws.Range("D9").Value = "2.604.15"
$ws.Range("E9").Value = "  -0.75%  "

$ws.Range("E10").Value = "  -1.88%  "

$ws.Range("E11").Value = "  +0.21%  "

$ws.Range("D12").Value = "'0.155"
$ws.Range("E12").Value = "  -3.28%  "

$ws.Range("E13").Value = "  -0.77%  "

$ws.Range("D14").Value = "3.058.53"
$ws.Range("E14").Value = "  -0.26%  "

$ws.Range("D15").Value = "'24.24"
$ws.Range("E15").Value = "  +3.60%  "

$ws.Range("D16").Value = "60.134.73"
$ws.Range("E16").Value = "  +0.65%  "

$ws.Range("D17").Value = "'0.0000140"
$ws.Range("E17").Value = "  +1.70%  "

$ws.Range("D18").Value = "2.603.10"
$ws.Range("E18").Value = "  -0.40%  "

$ws.Range("D19").Value = "'11.32"
$ws.Range("E19").Value = "  +6.10%  "

$ws.Range("E20").Value = "  -0.37%  "

$ws.Range("D21").Value = "'345.66"
$ws.Range("E21").Value = "  +0.78%  "

$ws.Range("E22").Value = "  +0.74%  "

$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("E24").Value = "  +1.57%  "

$ws.Range("D25").Value = "'63.03"
$ws.Range("E25").Value = "  +0.79%  "

$ws.Range("E26").Value = "  +0.33%  "

$ws.Range("E27").Value = "  -0.46%  "

$ws.Range("D28").Value = "'7.99"
$ws.Range("E28").Value = "  +4.41%  "

$ws.Range("D29").Value = "0.0₃0793"
$ws.Range("E29").Value = "  +1.44%  "

$ws.Range("E30").Value = "  +10.06%  "

$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  +0.09%  "

$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'6.36"
$ws.Range("E32").Value = "  +3.07%  "

$ws.Range("D33").Value = "'166.28"
$ws.Range("E33").Value = "  +4.91%  "

$ws.Range("D34").Value = "'19.40"
$ws.Range("E34").Value = "  -0.29%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'4.26"
$ws.Range("E35").Value = "  +2.51%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.29"
$ws.Range("E36").Value = "  +8.69%  "

$ws.Range("D37").Value = "'0.980"
$ws.Range("E37").Value = "  +5.02%  "

$ws.Range("E38").Value = "  +6.18%  "

$ws.Range("E39").Value = "  +0.74%  "

$ws.Range("D40").Value = "'311.34"
$ws.Range("E40").Value = "  +5.60%  "

$ws.Range("D41").Value = "'3.88"
$ws.Range("E41").Value = "  +4.20%  "

$ws.Range("D42").Value = "'0.836"
$ws.Range("E42").Value = "  -1.21%  "

$ws.Range("D43").Value = "'135.65"
$ws.Range("E43").Value = "  -3.82%  "

$ws.Range("D44").Value = "'0.0993"
$ws.Range("E44").Value = "  +1.47%  "

$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.28%  "

$ws.Range("D46").Value = "'19.84"
$ws.Range("E46").Value = "  +2.20%  "

$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0550"
$ws.Range("E47").Value = "  +2.28%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.604"
$ws.Range("E48").Value = "  +0.24%  "

$ws.Range("D49").Value = "'4.95"
$ws.Range("E49").Value = "  +3.99%  "

$ws.Range("D50").Value = "'0.0240"
$ws.Range("E50").Value = "  +0.03%  "

$ws.Range("D51").Value = "'19.89"
$ws.Range("E51").Value = "  +4.27%  "

